$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Test Exp 13"
$ws.Range("C14").Value = 30
$ws.Range("D14").Value = 0.3
$ws.Range("E14").Value = "Local"
$ws.Range("F14").Value = -1
$ws.Range("G14").Value = "28*28"
$ws.Range("H14").Value = "32*32"
$ws.Range("I14").Value = "3,4,5"

$ws.Range("A14:H14").HorizontalAlignment = -4131

$ws.Range("E18").Select()
